# BOM_ZumoComSystem: add 3V3 & 5V Regulator rows, add discharge-path rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row above the old "Taster" row (row 12) so the existing
# "3,3V-Regler" row (row 11) keeps its place and a brand-new "5V-Regler"
# row is inserted right after it.
$ws.Range("A12:H12").EntireRow.Insert()

# Update the (still in place) 3,3V-Regler row with the new part + datasheet link
$ws.Range("B11").Value = "TPS62056DGS"
$ws.Range("F11").Value = "https://www.digikey.de/product-detail/de/texas-instruments/TPS62056DGS/296-14212-5-ND/526047"

# Fill in the newly-inserted 5V-Regler row
$ws.Range("A12").Value = "5V-Regler"
$ws.Range("B12").Value = "LT1374CS8-5#PBF"
$ws.Range("C12").Value = 1
$ws.Range("F12").Value = "https://www.digikey.de/product-detail/de/analog-devices-inc/LT1374CS8-5-PBF/LT1374CS8-5-PBF-ND/888771"

# Append the two new "discharge path" rows at the bottom of the BOM
$ws.Range("A18").Value = "Spule 3V3"
$ws.Range("A19").Value = "Spule 5V"

# Match the author's final selection
$ws.Range("B15").Select() | Out-Null
